$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.695.43'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.600.12'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'210.91"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('E6').Value = '  +1.63%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = "'0.0618"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').Value = "'19.64"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '1.821.69'
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('D13').Value = '1.589.05'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').Value = "'0.521"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = "'64.82"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').Value = '26.671.84'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').Value = "'209.01"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = "'6.77"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('E23').Value = '  -3.67%  '
$ws.Range('D24').Value = "'8.91"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = "'145.58"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -2.84%  '
$ws.Range('E28').Value = '  +2.24%  '
$ws.Range('D29').Value = "'15.30"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('D33').Value = "'0.660"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').Value = '1.295.75'
$ws.Range('E35').Value = '  -1.83%  '
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').Value = "'1.49"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').Value = "'0.846"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('D41').Value = "'5.41"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.25%  '
$ws.Range('E42').Value = '  +1.30%  '
$ws.Range('D43').Value = "'0.787"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').Value = "'63.74"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.95%  '
$ws.Range('D45').Value = '1.734.82'
$ws.Range('E45').Value = '  +0.93%  '
$ws.Range('D46').Value = "'0.900"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.93%  '
$ws.Range('D47').Value = "'90.05"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('D49').Value = "'0.101"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.73%  '
$ws.Range('D50').Value = "'0.0506"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₇0981'
$ws.Range('E51').Value = '  -0.01%  '
